$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-unused blank placeholder cells in the header rows ---
$ws.Range("A1").Clear()
$ws.Range("B1").Clear()
$ws.Range("D1").Clear()
$ws.Range("F1").Clear()
$ws.Range("A2").Clear()

# --- Update column headers (row 2): F/P -> Chisq/Pr(>Chisq) ---
$ws.Range("C2").Value = "Chisq"
$ws.Range("D2").Value = "Pr(>Chisq)"
$ws.Range("E2").Value = "Chisq"
$ws.Range("F2").Value = "Pr(>Chisq)"

# --- Update the stats table body with new model results ---
$ws.Range("C3").Value = 5.73724017104805
$ws.Range("D3").Value = 0.0166088698482464
$ws.Range("E3").Value = 1.46613684530701
$ws.Range("F3").Value = 0.225956202970915

$ws.Range("C4").Value = 0.0214269073555852
$ws.Range("D4").Value = 0.883621946228612
$ws.Range("E4").Value = 0.0150493769051595
$ws.Range("F4").Value = 0.902363747208357

$ws.Range("C5").Value = 0.376889869665604
$ws.Range("D5").Value = 0.539272443454048
$ws.Range("E5").Value = 1.6416222620774
$ws.Range("F5").Value = 0.200103032572333

# --- Remove the Residuals row (old row 6); table now ends at row 5 ---
$ws.Rows.Item(6).Delete()

# --- Re-fit column widths to the new (wider, more precise) values ---
$ws.Columns.Item(1).ColumnWidth = 11.917
$ws.Columns.Item(2).ColumnWidth = 2.251
$ws.Range("C1:F1").EntireColumn.ColumnWidth = 11.251

# --- Selection moves to the refreshed results block ---
$ws.Range("C3:F5").Select()
